# Optuna Attempt (go back with original)
#
# Reverts a handful of forecast-derived metrics on "Forecast Comparison"
# (Inventory Coverage / Seasonality Index columns) and the two dependent
# roll-up totals on "Summary" back to their prior values.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------

# Inventory Coverage (column H)
$wsForecast.Range("H2").Value = 9.199999999999999
$wsForecast.Range("H7").Value = 3.95

# Seasonality Index (column L)
$wsForecast.Range("L2").Value  = 1.04
$wsForecast.Range("L3").Value  = 0.86
$wsForecast.Range("L4").Value  = 0.88
$wsForecast.Range("L5").Value  = 1.16
$wsForecast.Range("L6").Value  = 1.04
$wsForecast.Range("L7").Value  = 0.89
$wsForecast.Range("L8").Value  = 0.93
$wsForecast.Range("L9").Value  = 1.08
$wsForecast.Range("L10").Value = 0.86
$wsForecast.Range("L11").Value = 1.03
$wsForecast.Range("L12").Value = 1.16
$wsForecast.Range("L13").Value = 1.07
$wsForecast.Range("L14").Value = 0.84
$wsForecast.Range("L15").Value = 1.12
$wsForecast.Range("L16").Value = 0.87
$wsForecast.Range("L17").Value = 0.82

# --- Summary sheet ---------------------------------------------------------
#
# B9/B11 hold digit-only text (originally written as inline strings, e.g.
# "1979"), not numbers. A plain Range.Value = "1978" assignment gets
# auto-coerced to a numeric cell by Excel (same as typing it by hand), and
# the usual work-arounds for forcing text (a leading apostrophe, or setting
# NumberFormat to "@") both stick a "quote prefix"/text number-format onto
# the cell's style, which the source file never had.
#
# Instead, stage the text value with a leading apostrophe in an unused
# scratch cell, copy it, and PasteSpecial *values only* into the target
# cell - this carries over the text type without carrying over the
# scratch cell's style. The scratch cell is fully cleared afterwards so it
# leaves no trace (no stray formatting, no dimension growth).

$scratch = $wsSummary.Range("Z100")

$scratch.Value = "'1978"
$scratch.Copy()
$wsSummary.Range("B9").PasteSpecial(-4163)

$scratch.Value = "'522"
$scratch.Copy()
$wsSummary.Range("B11").PasteSpecial(-4163)

$scratch.Clear()
